$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# Row 2
Set-TextValue "D2" "288.78"
Set-TextValue "E2" "-0.59%"
Set-TextValue "G2" "6"

# Row 3
Set-TextValue "D3" "30.99"
Set-TextValue "E3" "0.47%"
Set-TextValue "G3" "6"

# Row 4
Set-TextValue "D4" "4.919"
Set-TextValue "E4" "-0.30%"
Set-TextValue "G4" "6"

# Row 5
Set-TextValue "D5" "0.07397"
Set-TextValue "E5" "2.66%"
Set-TextValue "G5" "6"

# Row 6
Set-TextValue "D6" "2.239"
Set-TextValue "E6" "26.13%"
Set-TextValue "G6" "6"

# Row 7
Set-TextValue "D7" "7.704"
Set-TextValue "E7" "0.51%"
Set-TextValue "G7" "6"

# Row 8
Set-TextValue "D8" "3.749"
Set-TextValue "E8" "0.66%"
Set-TextValue "G8" "6"

# Row 9
Set-TextValue "D9" "0.9101"
Set-TextValue "E9" "1.60%"
Set-TextValue "G9" "6"

# Row 10
Set-TextValue "D10" "0.08789"
Set-TextValue "E10" "13.80%"
Set-TextValue "G10" "6"

# Row 11
Set-TextValue "D11" "0.1690"
Set-TextValue "E11" "1.67%"
Set-TextValue "G11" "6"

# Row 12
Set-TextValue "D12" "0.08282"
Set-TextValue "E12" "3.57%"
Set-TextValue "G12" "6"

# Row 13
Set-TextValue "D13" "0.03112"
Set-TextValue "E13" "2.71%"
Set-TextValue "G13" "6"

# Row 14
Set-TextValue "D14" "0.09956"
Set-TextValue "E14" "-0.70%"
Set-TextValue "G14" "6"

# Row 15
Set-TextValue "D15" "0.001503"
Set-TextValue "E15" "-0.19%"
Set-TextValue "G15" "6"

# Row 16
Set-TextValue "D16" "0.005838"
Set-TextValue "E16" "0.80%"
Set-TextValue "G16" "6"

# Row 17
Set-TextValue "D17" "3.491"
Set-TextValue "E17" "0.43%"
Set-TextValue "G17" "6"

# Row 18
Set-TextValue "G18" "6"

# Row 19
Set-TextValue "D19" "0.3329"
Set-TextValue "E19" "1.55%"
Set-TextValue "G19" "6"

# Row 20
Set-TextValue "D20" "0.1297"
Set-TextValue "E20" "-1.27%"
Set-TextValue "G20" "6"

# Row 21
Set-TextValue "D21" "3.845"
Set-TextValue "E21" "-4.80%"
Set-TextValue "G21" "6"

# Row 22
Set-TextValue "G22" "6"

# Row 23
Set-TextValue "D23" "0.04564"
Set-TextValue "E23" "0.98%"
Set-TextValue "G23" "6"

# Row 24
Set-TextValue "D24" "0.001210"
Set-TextValue "E24" "-0.28%"
Set-TextValue "G24" "6"

# Row 25
Set-TextValue "D25" "0.004577"
Set-TextValue "E25" "14.19%"
Set-TextValue "G25" "6"

# Row 26
Set-TextValue "D26" "0.0001303"
Set-TextValue "E26" "4.16%"
Set-TextValue "G26" "6"

# Row 27
Set-TextValue "D27" "0.0003399"
Set-TextValue "E27" "-95.48%"
Set-TextValue "G27" "6"

# Row 28
Set-TextValue "G28" "6"

# Row 29
Set-TextValue "G29" "6"

# Row 30
Set-TextValue "G30" "6"

# Row 31
Set-TextValue "G31" "6"

# Row 32
Set-TextValue "G32" "6"

# Row 33
Set-TextValue "G33" "6"

# Row 34
Set-TextValue "G34" "6"

# Row 35
Set-TextValue "G35" "6"

# Row 36
Set-TextValue "G36" "6"

# Row 37
Set-TextValue "G37" "6"

# Row 38
Set-TextValue "G38" "6"

# Row 39
Set-TextValue "D39" "0.01585"
Set-TextValue "E39" "-0.33%"
Set-TextValue "G39" "6"

# Row 40
Set-TextValue "D40" "0.04466"
Set-TextValue "E40" "1.85%"
Set-TextValue "G40" "6"

# Row 41
Set-TextValue "D41" "0.007343"
Set-TextValue "E41" "-0.66%"
Set-TextValue "G41" "6"

# Row 42
Set-TextValue "D42" "0.009592"
Set-TextValue "E42" "25.00%"
Set-TextValue "G42" "6"

# Row 43
Set-TextValue "D43" "0.1324"
Set-TextValue "E43" "1.25%"
Set-TextValue "G43" "6"

# Row 44
Set-TextValue "D44" "0.002235"
Set-TextValue "E44" "8.95%"
Set-TextValue "G44" "6"

# Row 45
Set-TextValue "D45" "0.008876"
Set-TextValue "E45" "-4.42%"
Set-TextValue "G45" "6"

# Row 46
Set-TextValue "D46" "0.00006098"
Set-TextValue "E46" "2.75%"
Set-TextValue "G46" "6"

# Row 47
Set-TextValue "E47" "0.14%"
Set-TextValue "G47" "6"

# Row 48
Set-TextValue "D48" "2.197"
Set-TextValue "E48" "-2.14%"
Set-TextValue "G48" "6"

# Row 49
Set-TextValue "D49" "0.002003"
Set-TextValue "E49" "-33.27%"
Set-TextValue "G49" "6"

# Row 50
Set-TextValue "D50" "0.00002104"
Set-TextValue "E50" "0.14%"
Set-TextValue "G50" "6"

# Row 51
Set-TextValue "D51" "0.0002004"
Set-TextValue "E51" "0.14%"
Set-TextValue "G51" "6"
